# Update Csf3-Csf3r NATMI output with the new TPM-derived numbers.
# The original sheet had 3 data rows (sending cluster fixed at "MuSCs",
# target cluster cycling through ECs / MuSCs / Resolving-Mac). The new
# data has 6 rows: sending cluster cycles through ECs / MuSCs /
# Resolving-Mac, target cluster cycles through ECs / Resolving-Mac
# (MuSCs is no longer a receiving cluster for Csf3r under the new TPM).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (2-4) before laying down the new table
# (rows 2-7) so no stale cells are left behind.
$ws.Range("A2:T7").ClearContents()

$rows = @(
    @("ECs",          "Csf3", "Csf3r", "ECs",           2, 0.6666666666666666, 0.6262416666666667, 1.878725,           0.664168335284768,  0.6641683352847682, 2, 0.6666666666666666, 0.054819,           0.164457,  0.001660638168858053, 0.001660638168858053, 0.034329941925,     0.308969477325,     0.001102943288120799, 0.001102943288120799),
    @("ECs",          "Csf3", "Csf3r", "Resolving-Mac", 2, 0.6666666666666666, 0.6262416666666667, 1.878725,           0.664168335284768,  0.6641683352847682, 3, 1,                  32.95598433333333, 98.867953, 0.9983393618311419,   0.9983393618311419,   20.63841055554722,  185.745694999925,  0.6630653919966473,   0.6630653919966474),
    @("MuSCs",        "Csf3", "Csf3r", "ECs",           1, 0.3333333333333333, 0.083955,           0.251865,           0.08903951231100779,0.08903951231100779,2, 0.6666666666666666, 0.054819,           0.164457,  0.001660638168858053, 0.001660638168858053, 0.004602329145,     0.041420962305,     0.0001478624126801661,0.0001478624126801661),
    @("MuSCs",        "Csf3", "Csf3r", "Resolving-Mac", 1, 0.3333333333333333, 0.083955,           0.251865,           0.08903951231100779,0.08903951231100779,3, 1,                  32.95598433333333, 98.867953, 0.9983393618311419,   0.9983393618311419,   2.766819664705,     24.901376982345,    0.08889164989832762,  0.08889164989832762),
    @("Resolving-Mac","Csf3", "Csf3r", "ECs",           2, 0.6666666666666666, 0.2326993333333333, 0.698098,           0.2467921524042241, 0.2467921524042241, 2, 0.6666666666666666, 0.054819,           0.164457,  0.001660638168858053, 0.001660638168858053, 0.012756344754,     0.114807102786,     0.0004098324680570884,0.0004098324680570884),
    @("Resolving-Mac","Csf3", "Csf3r", "Resolving-Mac", 2, 0.6666666666666666, 0.2326993333333333, 0.698098,           0.2467921524042241, 0.2467921524042241, 3, 1,                  32.95598433333333, 98.867953, 0.9983393618311419,   0.9983393618311419,   7.668835583710445,  69.01952025339401,  0.246382319936167,    0.246382319936167)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
